# Apply the "concise marksheet" corrections to the quiz marksheet.
# Commit message: "changes in concise marksheet / Corr/total marks"
# The Total row's correct-answer count and the Correct/Total marks summary
# text are updated to reflect 110 correct out of 140 total marks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: B11 (Right marking count) 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row: B12 (Right total) 66 -> 110
$ws.Range("B12").Value = 110

# "Total" row: E12 (Correct/Total marks text) "63/84" -> "110/140"
$ws.Range("E12").Value = "110/140"
